$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new header cells (AD1:AF1) ---
# Copy formatting (bold font, border, alignment) from the last existing
# header cell (AC1) onto the new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Fill in the season record (Wins/Losses/Ties) for every player row ---
$wins = 96
$losses = 67
$ties = 0

$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
